$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: volume/number and date range
$ws.Range("A8").Value = "Volume 31   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/8/2024  Through  4/14/2024"

# Data table updates (rows 15-31)
$ws.Range("C15").Value = 1
$ws.Range("G15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("F15").Value = 1
$ws.Range("G15").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -25
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 142.857142857143
$ws.Range("I16").Value = 46
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 15
$ws.Range("L16").Value = -9.803921568627
$ws.Range("M16").Value = -19.298245614035
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 6.25
$ws.Range("I17").Value = 72
$ws.Range("J17").Value = 71
$ws.Range("K17").Value = 1.408450704225
$ws.Range("L17").Value = -13.253012048192
$ws.Range("M17").Value = 46.938775510204
$ws.Range("C18").Value = 2
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 133.333333333333
$ws.Range("I18").Value = 36
$ws.Range("K18").Value = -21.739130434782
$ws.Range("L18").Value = 24.137931034482
$ws.Range("M18").Value = 33.333333333333
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 21.428571428571
$ws.Range("I19").Value = 101
$ws.Range("J19").Value = 84
$ws.Range("K19").Value = 20.238095238095
$ws.Range("L19").Value = 8.602150537634
$ws.Range("M19").Value = 48.529411764705
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = "'0"
$ws.Range("A20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").Value = "***.*"
$ws.Range("A20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 26
$ws.Range("K20").Value = -29.729729729729
$ws.Range("L20").Value = -39.534883720930
$ws.Range("M20").Value = 13.043478260869
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 12.5
$ws.Range("F21").Value = 82
$ws.Range("H21").Value = 34.426229508196
$ws.Range("I21").Value = 284
$ws.Range("J21").Value = 282
$ws.Range("K21").Value = 0.709219858156
$ws.Range("L21").Value = -6.270627062706
$ws.Range("M21").Value = 23.478260869565
$ws.Range("D22").Value = 1
$ws.Range("I22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("K22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("I22").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = -100
$ws.Range("K22").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -72.727272727272
$ws.Range("D23").Value = "'0"
$ws.Range("A23").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = "***.*"
$ws.Range("A23").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("L23").Value = -33.333333333333
$ws.Range("M23").Value = 100
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = -24
$ws.Range("F24").Value = 71
$ws.Range("G24").Value = 79
$ws.Range("H24").Value = -10.126582278481
$ws.Range("I24").Value = 230
$ws.Range("J24").Value = 272
$ws.Range("K24").Value = -15.441176470588
$ws.Range("L24").Value = -50.959488272921
$ws.Range("M24").Value = 76.923076923076
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 5
$ws.Range("I25").Value = 47
$ws.Range("J25").Value = 72
$ws.Range("K25").Value = -34.722222222222
$ws.Range("L25").Value = -84.228187919463
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = -38.888888888888
$ws.Range("I26").Value = 92
$ws.Range("J26").Value = 117
$ws.Range("K26").Value = -21.367521367521
$ws.Range("L26").Value = -12.380952380952
$ws.Range("M26").Value = -20
$ws.Range("C27").Value = 1
$ws.Range("G27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("A27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$ws.Range("A27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 1
$ws.Range("G27").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 4
$ws.Range("K27").Value = -20
$ws.Range("L27").Value = -20
$ws.Range("C28").Value = 1
$ws.Range("D28").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 12
$ws.Range("J28").Value = 11
$ws.Range("K28").Value = 9.090909090909
$ws.Range("L28").Value = -7.692307692307
$ws.Range("D31").Value = "'0"
$ws.Range("A31").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "***.*"
$ws.Range("A31").Copy()
$ws.Range("E31").PasteSpecial(-4122)
